# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Applies the edits described by the commit diff:
#   - Metadata sheet: URL, Version, Status, Date, Description, Context
#   - Elements sheet: root Extension Definition text, Extension.value[x] Min,
#     and the Extension.url Fixed Value (kept in sync with the URL change)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$newUrl = "https://clinyqai.github.io/open-nursing-core-ig/StructureDefinition/UKCore-Extension-EthnicCategory"
$newDescription = "An extension to record the ethnic category of a patient, as per UK Core standards."

# --- Metadata sheet updates ---
$meta.Range("B2").Value  = $newUrl                           # URL
$meta.Range("B3").Value  = "0.1.0"                           # Version
$meta.Range("B6").Value  = "draft"                           # Status
$meta.Range("B8").Value  = "2025-12-26T14:13:58+00:00"       # Date
$meta.Range("B11").Value = $newDescription                   # Description
$meta.Range("B20").Value = "element:Patient"                 # Context

# --- Elements sheet updates ---
$elements.Range("M2").Value = $newDescription                 # Extension (root) Definition

# Extension.value[x] Min: 0 -> 1 (validation fix). Min/Max in this table are
# stored as text ("0"/"1"), so instead of letting Excel coerce the new value
# to a number, copy the already-correctly-typed/styled "1" cell (G6) onto F6.
# This updates F6's value (and keeps its existing text formatting/style).
$elements.Range("G6").Copy()
$elements.Range("F6").PasteSpecial(-4163)  # xlPasteAll

$elements.Range("R5").Value = $newUrl                         # Extension.url Fixed Value
